$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F updates
$changes = @{
    "展览" = @(
        @{Row = 4;  Value = 1977},
        @{Row = 5;  Value = 5721},
        @{Row = 19; Value = 51},
        @{Row = 20; Value = 173},
        @{Row = 30; Value = 407},
        @{Row = 34; Value = 585},
        @{Row = 36; Value = 1737},
        @{Row = 37; Value = 2248},
        @{Row = 41; Value = 633},
        @{Row = 42; Value = 365}
    );
    "演出" = @(
        @{Row = 23; Value = 4}
    );
    "全部类型" = @(
        @{Row = 5;  Value = 5721},
        @{Row = 19; Value = 51},
        @{Row = 20; Value = 173},
        @{Row = 31; Value = 407},
        @{Row = 33; Value = 586},
        @{Row = 35; Value = 1737},
        @{Row = 42; Value = 365},
        @{Row = 47; Value = 4}
    )
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $changes[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
